$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:48:48"
$ws1.Range("A3").Value = "Total filas: 16"

# Existing rows 7-12: update Hora_Scrap (A) and Minutos (D)
$ws1.Cells.Item(7,1).Value = "04:48:48"
$ws1.Cells.Item(7,4).Value = 5

$ws1.Cells.Item(8,1).Value = "04:48:48"
$ws1.Cells.Item(8,4).Value = 28

$ws1.Cells.Item(9,1).Value = "04:48:48"
$ws1.Cells.Item(9,4).Value = 34

$ws1.Cells.Item(10,1).Value = "04:48:48"
$ws1.Cells.Item(10,4).Value = 46

$ws1.Cells.Item(11,1).Value = "04:48:48"
$ws1.Cells.Item(11,4).Value = 58

$ws1.Cells.Item(12,1).Value = "04:48:48"
$ws1.Cells.Item(12,4).Value = 66

# New rows 13-21
$ws1.Cells.Item(13,1).Value = "04:48:48"
$ws1.Cells.Item(13,2).Value = "06:04"
$ws1.Cells.Item(13,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(13,4).Value = 76
$ws1.Cells.Item(13,5).Value = "LP1912"

$ws1.Cells.Item(14,1).Value = "04:48:48"
$ws1.Cells.Item(14,2).Value = "06:11"
$ws1.Cells.Item(14,3).Value = "215A_EL PATO"
$ws1.Cells.Item(14,4).Value = 83
$ws1.Cells.Item(14,5).Value = "LP1912"

$ws1.Cells.Item(15,1).Value = "04:48:48"
$ws1.Cells.Item(15,2).Value = "06:14"
$ws1.Cells.Item(15,3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(15,4).Value = 86
$ws1.Cells.Item(15,5).Value = "LP1912"

$ws1.Cells.Item(16,1).Value = "04:48:48"
$ws1.Cells.Item(16,2).Value = "06:21"
$ws1.Cells.Item(16,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(16,4).Value = 93
$ws1.Cells.Item(16,5).Value = "LP1912"

$ws1.Cells.Item(17,1).Value = "04:48:48"
$ws1.Cells.Item(17,2).Value = "06:27"
$ws1.Cells.Item(17,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(17,4).Value = 99
$ws1.Cells.Item(17,5).Value = "LP1912"

$ws1.Cells.Item(18,1).Value = "04:48:48"
$ws1.Cells.Item(18,2).Value = "06:29"
$ws1.Cells.Item(18,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(18,4).Value = 101
$ws1.Cells.Item(18,5).Value = "LP1912"

$ws1.Cells.Item(19,1).Value = "04:48:48"
$ws1.Cells.Item(19,2).Value = "06:31"
$ws1.Cells.Item(19,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(19,4).Value = 103
$ws1.Cells.Item(19,5).Value = "LP1912"

$ws1.Cells.Item(20,1).Value = "04:48:48"
$ws1.Cells.Item(20,2).Value = "06:44"
$ws1.Cells.Item(20,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(20,4).Value = 116
$ws1.Cells.Item(20,5).Value = "LP1912"

$ws1.Cells.Item(21,1).Value = "04:48:48"
$ws1.Cells.Item(21,2).Value = "06:46"
$ws1.Cells.Item(21,3).Value = "215C_EL PATO"
$ws1.Cells.Item(21,4).Value = 118
$ws1.Cells.Item(21,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:48:48"
$ws2.Range("A3").Value = "Total filas: 4"

# Existing row 7: update Hora_Scrap (A) and Minutos (D)
$ws2.Cells.Item(7,1).Value = "04:48:48"
$ws2.Cells.Item(7,4).Value = 46

# New rows 8-9
$ws2.Cells.Item(8,1).Value = "04:48:48"
$ws2.Cells.Item(8,2).Value = "06:11"
$ws2.Cells.Item(8,3).Value = "215A_EL PATO"
$ws2.Cells.Item(8,4).Value = 83
$ws2.Cells.Item(8,5).Value = "LP1912"

$ws2.Cells.Item(9,1).Value = "04:48:48"
$ws2.Cells.Item(9,2).Value = "06:46"
$ws2.Cells.Item(9,3).Value = "215C_EL PATO"
$ws2.Cells.Item(9,4).Value = 118
$ws2.Cells.Item(9,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:48:48"
$ws3.Range("A3").Value = "Total filas: 4"

# Insert a new row above the existing data row (old row 6 becomes row 7)
$ws3.Rows.Item(6).Insert()

$ws3.Cells.Item(6,1).Value = "04:48:48"
$ws3.Cells.Item(6,2).Value = "05:43"
$ws3.Cells.Item(6,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(6,4).Value = 55
$ws3.Cells.Item(6,5).Value = "L6173"

# New rows 8-9 (appended after the shifted original row 7)
$ws3.Cells.Item(8,1).Value = "04:48:48"
$ws3.Cells.Item(8,2).Value = "06:08"
$ws3.Cells.Item(8,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(8,4).Value = 80
$ws3.Cells.Item(8,5).Value = "L6173"

$ws3.Cells.Item(9,1).Value = "04:48:48"
$ws3.Cells.Item(9,2).Value = "06:32"
$ws3.Cells.Item(9,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(9,4).Value = 104
$ws3.Cells.Item(9,5).Value = "L6203"
